# Applies the CryCompanywiseStockReport_1 stock-report corrections:
#  - item rows that were re-sorted (their code/desc/price/qty/amount swap
#    with the next row), and quantity/amount corrections on other rows.
#  - Sub Total / Grand Total cells are updated to the recomputed sums.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 180
$ws.Range("B180").Value = 57756
$ws.Range("F180").Value = 142
$ws.Range("G180").Value = 9434.48

# Row 181
$ws.Range("B181").Value = 53925
$ws.Range("F181").Value = 1
$ws.Range("G181").Value = 66.44

# Row 183
$ws.Range("F183").Value = 334
$ws.Range("G183").Value = 32334.54

# Row 186
$ws.Range("B186").Value = 80865.67999999999

# Row 193
$ws.Range("F193").Value = 9
$ws.Range("G193").Value = 14780.07

# Row 200
$ws.Range("B200").Value = 99054.09

# Row 256
$ws.Range("F256").Value = 1763
$ws.Range("G256").Value = 32615.5

# Row 261
$ws.Range("B261").Value = 43559.04

# Row 336
$ws.Range("B336").Value = 57854
$ws.Range("F336").Value = 2
$ws.Range("G336").Value = 611.6799999999999

# Row 337
$ws.Range("B337").Value = 62997
$ws.Range("F337").Value = 116
$ws.Range("G337").Value = 35477.44

# Row 339
$ws.Range("B339").Value = 61610
$ws.Range("D339").Value = 102.71
$ws.Range("E339").Value = 122.71
$ws.Range("F339").Value = 62
$ws.Range("G339").Value = 6368.02

# Row 340
$ws.Range("B340").Value = 57077
$ws.Range("D340").Value = 93.08
$ws.Range("E340").Value = 111.2
$ws.Range("F340").Value = 1
$ws.Range("G340").Value = 93.08

# Row 348
$ws.Range("F348").Value = 172
$ws.Range("G348").Value = 15201.36

# Row 351
$ws.Range("F351").Value = 57
$ws.Range("G351").Value = 9046.469999999999

# Row 382
$ws.Range("F382").Value = 201
$ws.Range("G382").Value = 11852.97

# Row 392
$ws.Range("B392").Value = 57870
$ws.Range("F392").Value = 0
$ws.Range("G392").Value = 0

# Row 393
$ws.Range("B393").Value = 63040
$ws.Range("F393").Value = 68
$ws.Range("G393").Value = 7467.76

# Row 398
$ws.Range("B398").Value = 63112
$ws.Range("F398").Value = 287
$ws.Range("G398").Value = 14961.31

# Row 399
$ws.Range("B399").Value = 57885
$ws.Range("F399").Value = 4
$ws.Range("G399").Value = 208.52

# Row 401
$ws.Range("F401").Value = 38
$ws.Range("G401").Value = 4254.48

# Row 416
$ws.Range("F416").Value = 1355
$ws.Range("G416").Value = 16178.7

# Row 417
$ws.Range("B417").Value = 57817
$ws.Range("F417").Value = 3
$ws.Range("G417").Value = 239.43

# Row 418
$ws.Range("B418").Value = 62865
$ws.Range("F418").Value = 239
$ws.Range("G418").Value = 19074.59

# Row 421
$ws.Range("B421").Value = 53060
$ws.Range("C421").Value = 'HUL-REXONA COCONUT&amp;OLIVE OILS 4x100g'
$ws.Range("D421").Value = 109.82
$ws.Range("E421").Value = 131.19
$ws.Range("F421").Value = 1
$ws.Range("G421").Value = 109.82

# Row 422
$ws.Range("B422").Value = 63043
$ws.Range("C422").Value = 'HUL-Rexona Coconut&amp;Olive Oils 4X100G'
$ws.Range("D422").Value = 115.01
$ws.Range("E422").Value = 137.41
$ws.Range("F422").Value = 55
$ws.Range("G422").Value = 6325.55

# Row 428
$ws.Range("B428").Value = 62933
$ws.Range("F428").Value = 200
$ws.Range("G428").Value = 11826

# Row 429
$ws.Range("B429").Value = 57835
$ws.Range("F429").Value = 1
$ws.Range("G429").Value = 59.13

# Row 436
$ws.Range("F436").Value = 5629
$ws.Range("G436").Value = 118377.87

# Row 438
$ws.Range("B438").Value = 62784
$ws.Range("F438").Value = 14
$ws.Range("G438").Value = 1163.4

# Row 439
$ws.Range("B439").Value = 57799
$ws.Range("F439").Value = 2
$ws.Range("G439").Value = 166.2

# Row 440
$ws.Range("F440").Value = 165
$ws.Range("G440").Value = 26565

# Row 448
$ws.Range("B448").Value = 57856
$ws.Range("F448").Value = 2
$ws.Range("G448").Value = 342.66

# Row 449
$ws.Range("B449").Value = 63007
$ws.Range("F449").Value = 1109
$ws.Range("G449").Value = 190004.97

# Row 450
$ws.Range("B450").Value = 57857
$ws.Range("F450").Value = 3
$ws.Range("G450").Value = 453.51

# Row 451
$ws.Range("B451").Value = 63008
$ws.Range("F451").Value = 615
$ws.Range("G451").Value = 92969.55

# Row 460
$ws.Range("F460").Value = 654
$ws.Range("G460").Value = 8861.700000000001

# Row 461
$ws.Range("B461").Value = 53082
$ws.Range("C461").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F461").Value = 1
$ws.Range("G461").Value = 59.47

# Row 462
$ws.Range("B462").Value = 63102
$ws.Range("C462").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F462").Value = 199
$ws.Range("G462").Value = 11834.53

# Row 463
$ws.Range("F463").Value = 2
$ws.Range("G463").Value = 595.28

# Row 464
$ws.Range("B464").Value = 1345809.27

# Row 466
$ws.Range("F466").Value = 53
$ws.Range("G466").Value = 9729.74

# Row 481
$ws.Range("B481").Value = 87562.99000000001

# Row 490
$ws.Range("F490").Value = 180
$ws.Range("G490").Value = 5803.2

# Row 491
$ws.Range("F491").Value = 90
$ws.Range("G491").Value = 2017.8

# Row 496
$ws.Range("B496").Value = 27114.01

# Row 610
$ws.Range("F610").Value = 209
$ws.Range("G610").Value = 3925.02

# Row 627
$ws.Range("B627").Value = 101647.72

# Row 646
$ws.Range("F646").Value = 3
$ws.Range("G646").Value = 483.6

# Row 649
$ws.Range("B649").Value = 19973.85

# Row 682
$ws.Range("F682").Value = 45
$ws.Range("G682").Value = 1306.35

# Row 688
$ws.Range("F688").Value = 119
$ws.Range("G688").Value = 4063.85

# Row 696
$ws.Range("B696").Value = 48129.5

# Row 733
$ws.Range("F733").Value = 89
$ws.Range("G733").Value = 4884.32

# Row 739
$ws.Range("F739").Value = 37
$ws.Range("G739").Value = 2349.13

# Row 746
$ws.Range("B746").Value = 86666.5

# Row 797
$ws.Range("F797").Value = 32
$ws.Range("G797").Value = 3899.2

# Row 804
$ws.Range("B804").Value = 84360.67

# Row 967
$ws.Range("B967").Value = 5708894.88

# Row 968
$ws.Range("B968").Value = 5708894.88
